# Auto-generated Excel COM-interop script
# Applies updated market-board price snapshot values (columns H:N)
# to each class sheet's Leve profit table, per the commit diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 478.07693
$ws.Range("I33").Value = 492.27274
$ws.Range("J33").Value = 400
$ws.Range("K33").Value = 492.27274
$ws.Range("L33").Value = 400
$ws.Range("M33").Value = -263.27274
$ws.Range("N33").Value = -858
$ws.Range("H62").Value = 3363.4
$ws.Range("I62").Value = 1666.8
$ws.Range("J62").Value = 5060
$ws.Range("K62").Value = 1666.8
$ws.Range("L62").Value = 5060
$ws.Range("M62").Value = -1042.8
$ws.Range("N62").Value = -6308
$ws.Range("H65").Value = 3363.4
$ws.Range("I65").Value = 1666.8
$ws.Range("J65").Value = 5060
$ws.Range("K65").Value = 8334
$ws.Range("L65").Value = 25300
$ws.Range("M65").Value = -5214
$ws.Range("N65").Value = -31540
$ws.Range("H88").Value = 2837.4167
$ws.Range("I88").Value = 10000
$ws.Range("J88").Value = 2186.2727
$ws.Range("K88").Value = 10000
$ws.Range("L88").Value = 2186.2727
$ws.Range("M88").Value = -9594
$ws.Range("N88").Value = -2998.2727
$ws.Range("H91").Value = 2837.4167
$ws.Range("I91").Value = 10000
$ws.Range("J91").Value = 2186.2727
$ws.Range("K91").Value = 10000
$ws.Range("L91").Value = 2186.2727
$ws.Range("M91").Value = -8596
$ws.Range("N91").Value = -4994.2727
$ws.Range("H111").Value = 959.25
$ws.Range("I111").Value = 841.7778
$ws.Range("J111").Value = 1110.2858
$ws.Range("K111").Value = 2525.3334
$ws.Range("L111").Value = 3330.8574
$ws.Range("M111").Value = 541.6666
$ws.Range("N111").Value = -9464.857400000001
$ws.Range("H116").Value = 222540.66
$ws.Range("I116").Value = 418969.34
$ws.Range("J116").Value = 8254.817999999999
$ws.Range("K116").Value = 418969.34
$ws.Range("L116").Value = 8254.817999999999
$ws.Range("M116").Value = -415527.34
$ws.Range("N116").Value = -15138.818
$ws.Range("H118").Value = 893.7143
$ws.Range("J118").Value = 992.375
$ws.Range("L118").Value = 2977.125
$ws.Range("N118").Value = -6291.125
$ws.Range("H129").Value = 878.27
$ws.Range("J129").Value = 906.883
$ws.Range("L129").Value = 2720.649
$ws.Range("N129").Value = -12720.649
$ws.Range("H132").Value = 298645
$ws.Range("I132").Value = 4426.5806
$ws.Range("K132").Value = 13279.7418
$ws.Range("M132").Value = -10749.7418
$ws.Range("H136").Value = 49352
$ws.Range("J136").Value = 49352
$ws.Range("L136").Value = 49352
$ws.Range("N136").Value = -59552
$ws.Range("H137").Value = 2592.9512
$ws.Range("I137").Value = 1383.258
$ws.Range("J137").Value = 6343
$ws.Range("K137").Value = 4149.774
$ws.Range("L137").Value = 19029
$ws.Range("M137").Value = -1599.774
$ws.Range("N137").Value = -24129
$ws.Range("H138").Value = 2195.59
$ws.Range("J138").Value = 2448.7195
$ws.Range("L138").Value = 7346.1585
$ws.Range("N138").Value = -17626.1585

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4422.774
$ws.Range("I32").Value = 3561.446
$ws.Range("J32").Value = 7369.421
$ws.Range("K32").Value = 3561.446
$ws.Range("L32").Value = 7369.421
$ws.Range("M32").Value = -3274.446
$ws.Range("N32").Value = -7943.421
$ws.Range("H45").Value = 1391.3334
$ws.Range("I45").Value = 1530.4
$ws.Range("J45").Value = 1217.5
$ws.Range("K45").Value = 1530.4
$ws.Range("L45").Value = 1217.5
$ws.Range("M45").Value = -1153.4
$ws.Range("N45").Value = -1971.5
$ws.Range("H74").Value = 8982.182000000001
$ws.Range("I74").Value = 8380.4
$ws.Range("J74").Value = 15000
$ws.Range("K74").Value = 8380.4
$ws.Range("L74").Value = 15000
$ws.Range("M74").Value = -7506.4
$ws.Range("N74").Value = -16748
$ws.Range("H77").Value = 8982.182000000001
$ws.Range("I77").Value = 8380.4
$ws.Range("J77").Value = 15000
$ws.Range("K77").Value = 41902
$ws.Range("L77").Value = 75000
$ws.Range("M77").Value = -37534
$ws.Range("N77").Value = -83736
$ws.Range("H111").Value = 50644
$ws.Range("J111").Value = 50644
$ws.Range("L111").Value = 50644
$ws.Range("N111").Value = -58824
$ws.Range("H121").Value = 28855
$ws.Range("J121").Value = 28855
$ws.Range("L121").Value = 28855
$ws.Range("N121").Value = -32349
$ws.Range("H132").Value = 1910.6875
$ws.Range("I132").Value = 1174.7693
$ws.Range("J132").Value = 5099.6665
$ws.Range("K132").Value = 3524.3079
$ws.Range("L132").Value = 15298.9995
$ws.Range("M132").Value = -994.3078999999998
$ws.Range("N132").Value = -20358.9995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1208.6154
$ws.Range("I107").Value = 1208.6154
$ws.Range("K107").Value = 1208.6154
$ws.Range("M107").Value = 711.3846000000001
$ws.Range("H134").Value = 1823.2433
$ws.Range("I134").Value = 1122.2
$ws.Range("J134").Value = 3283.75
$ws.Range("K134").Value = 3366.6
$ws.Range("L134").Value = 9851.25
$ws.Range("M134").Value = -831.6000000000004
$ws.Range("N134").Value = -14921.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 20836398
$ws.Range("I31").Value = 1039.5714
$ws.Range("J31").Value = 50005900
$ws.Range("K31").Value = 1039.5714
$ws.Range("L31").Value = 50005900
$ws.Range("M31").Value = -744.5714
$ws.Range("N31").Value = -50006490
$ws.Range("H34").Value = 20836398
$ws.Range("I34").Value = 1039.5714
$ws.Range("J34").Value = 50005900
$ws.Range("K34").Value = 1039.5714
$ws.Range("L34").Value = 50005900
$ws.Range("M34").Value = -837.5714
$ws.Range("N34").Value = -50006304
$ws.Range("H94").Value = 1119.3125
$ws.Range("J94").Value = 1152
$ws.Range("L94").Value = 1152
$ws.Range("N94").Value = -2054
$ws.Range("H132").Value = 2659.121
$ws.Range("I132").Value = 1823.6
$ws.Range("J132").Value = 3944.5386
$ws.Range("K132").Value = 5470.799999999999
$ws.Range("L132").Value = 11833.6158
$ws.Range("M132").Value = -2940.799999999999
$ws.Range("N132").Value = -16893.6158
$ws.Range("H134").Value = 5739.76
$ws.Range("I134").Value = 6093.1577
$ws.Range("J134").Value = 4620.6665
$ws.Range("K134").Value = 18279.4731
$ws.Range("L134").Value = 13861.9995
$ws.Range("M134").Value = -15744.4731
$ws.Range("N134").Value = -18931.9995
$ws.Range("H140").Value = 74382.44
$ws.Range("J140").Value = 74382.44
$ws.Range("L140").Value = 74382.44
$ws.Range("N140").Value = -84742.44

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 5083.3335
$ws.Range("I63").Value = 4750
$ws.Range("K63").Value = 14250
$ws.Range("M63").Value = -13501
$ws.Range("H66").Value = 5083.3335
$ws.Range("I66").Value = 4750
$ws.Range("K66").Value = 42750
$ws.Range("M66").Value = -39006
$ws.Range("H113").Value = 639.17645
$ws.Range("I113").Value = 627.0769
$ws.Range("K113").Value = 1881.2307
$ws.Range("M113").Value = 288.7692999999999
$ws.Range("H122").Value = 2802.3096
$ws.Range("I122").Value = 1069.4
$ws.Range("K122").Value = 9624.6
$ws.Range("M122").Value = -7174.6
$ws.Range("H129").Value = 2867.3333
$ws.Range("J129").Value = 2723.6667
$ws.Range("L129").Value = 8171.000100000001
$ws.Range("N129").Value = -18171.0001
$ws.Range("H131").Value = 7353827
$ws.Range("I131").Value = 83333700
$ws.Range("J131").Value = 936.0323
$ws.Range("K131").Value = 250001100
$ws.Range("L131").Value = 2808.0969
$ws.Range("M131").Value = -249996060
$ws.Range("N131").Value = -12888.0969

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 20835916
$ws.Range("I80").Value = 35716340
$ws.Range("J80").Value = 3322.4
$ws.Range("K80").Value = 35716340
$ws.Range("L80").Value = 3322.4
$ws.Range("M80").Value = -35715342
$ws.Range("N80").Value = -5318.4
$ws.Range("H83").Value = 20835916
$ws.Range("I83").Value = 35716340
$ws.Range("J83").Value = 3322.4
$ws.Range("K83").Value = 178581700
$ws.Range("L83").Value = 16612
$ws.Range("M83").Value = -178576708
$ws.Range("N83").Value = -26596
$ws.Range("H132").Value = 2913.7896
$ws.Range("I132").Value = 1755.1666
$ws.Range("J132").Value = 4900
$ws.Range("K132").Value = 5265.4998
$ws.Range("L132").Value = 14700
$ws.Range("M132").Value = -2735.4998
$ws.Range("N132").Value = -19760

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 33824.97
$ws.Range("I22").Value = 67884.266
$ws.Range("J22").Value = 1894.375
$ws.Range("K22").Value = 67884.266
$ws.Range("L22").Value = 1894.375
$ws.Range("M22").Value = -67589.266
$ws.Range("N22").Value = -2484.375
$ws.Range("H27").Value = 33824.97
$ws.Range("I27").Value = 67884.266
$ws.Range("J27").Value = 1894.375
$ws.Range("K27").Value = 67884.266
$ws.Range("L27").Value = 1894.375
$ws.Range("M27").Value = -67777.266
$ws.Range("N27").Value = -2108.375
$ws.Range("H55").Value = 390.15384
$ws.Range("I55").Value = 272
$ws.Range("J55").Value = 491.42856
$ws.Range("K55").Value = 272
$ws.Range("L55").Value = 491.42856
$ws.Range("M55").Value = -99
$ws.Range("N55").Value = -837.4285600000001
$ws.Range("H118").Value = 28390
$ws.Range("J118").Value = 28390
$ws.Range("L118").Value = 28390
$ws.Range("N118").Value = -31704
$ws.Range("H132").Value = 17058.637
$ws.Range("I132").Value = 26299.2
$ws.Range("J132").Value = 9358.166999999999
$ws.Range("K132").Value = 78897.60000000001
$ws.Range("L132").Value = 28074.501
$ws.Range("M132").Value = -76367.60000000001
$ws.Range("N132").Value = -33134.501
$ws.Range("H136").Value = 3850.5715
$ws.Range("I136").Value = 1576
$ws.Range("J136").Value = 6883.3335
$ws.Range("K136").Value = 4728
$ws.Range("L136").Value = 20650.0005
$ws.Range("M136").Value = -2178
$ws.Range("N136").Value = -25750.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 18059.5
$ws.Range("J41").Value = 32777
$ws.Range("L41").Value = 32777
$ws.Range("N41").Value = -33557
$ws.Range("H81").Value = 25001198
$ws.Range("I81").Value = 31251298
$ws.Range("J81").Value = 798.5
$ws.Range("K81").Value = 62502596
$ws.Range("L81").Value = 1597
$ws.Range("M81").Value = -62501535
$ws.Range("N81").Value = -3719
$ws.Range("H84").Value = 25001198
$ws.Range("I84").Value = 31251298
$ws.Range("J84").Value = 798.5
$ws.Range("K84").Value = 312512980
$ws.Range("L84").Value = 7985
$ws.Range("M84").Value = -312507676
$ws.Range("N84").Value = -18593
$ws.Range("H132").Value = 8334853.5
$ws.Range("I132").Value = 929.41174
$ws.Range("J132").Value = 55560424
$ws.Range("K132").Value = 2788.23522
$ws.Range("L132").Value = 166681272
$ws.Range("M132").Value = -258.23522
$ws.Range("N132").Value = -166686332

Write-Host "Applied 278 cell updates across 8 sheets"
